$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2025-05-23 Friday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-05-24 Saturday", 2)

# Update the division expressions (each is unique in the document, so a simple
# Find/Replace is safe and order-independent)
$d.Content.Find.Execute("48÷5=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "81÷9=", 2)
$d.Content.Find.Execute("14÷2=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "61÷2=", 2)
$d.Content.Find.Execute("23÷4=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "92÷2=", 2)
$d.Content.Find.Execute("51÷3=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "55÷8=", 2)
$d.Content.Find.Execute("72÷4=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "30÷8=", 2)
$d.Content.Find.Execute("14÷7=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "79÷9=", 2)
$d.Content.Find.Execute("76÷8=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "74÷8=", 2)
$d.Content.Find.Execute("50÷4=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "14÷3=", 2)
$d.Content.Find.Execute("11÷2=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "73÷3=", 2)
$d.Content.Find.Execute("25÷9=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "62÷2=", 2)
$d.Content.Find.Execute("60÷7=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "95÷3=", 2)
$d.Content.Find.Execute("96÷2=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "37÷6=", 2)
$d.Content.Find.Execute("99÷6=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "86÷4=", 2)
$d.Content.Find.Execute("26÷9=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "58÷5=", 2)
$d.Content.Find.Execute("19÷8=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "48÷2=", 2)
$d.Content.Find.Execute("42÷8=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "76÷7=", 2)
$d.Content.Find.Execute("57÷2=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "80÷8=", 2)
$d.Content.Find.Execute("56÷3=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "31÷2=", 2)
$d.Content.Find.Execute("57÷4=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "59÷7=", 2)
$d.Content.Find.Execute("33÷9=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "81÷5=", 2)
$d.Content.Find.Execute("60÷2=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "77÷8=", 2)
$d.Content.Find.Execute("68÷8=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "64÷8=", 2)
$d.Content.Find.Execute("93÷4=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "33÷3=", 2)
$d.Content.Find.Execute("55÷5=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "20÷6=", 2)
$d.Content.Find.Execute("67÷9=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "98÷7=", 2)
